$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D1").Value = "Percentage"
$ws.Range("D1").Font.Bold = $true

# Percentage values (Number of Holders as % of the Male+Female total for
# that Qualification)
$ws.Range("D2").Value = 74.295907909900023
$ws.Range("D3").Value = 91.560186778955739
$ws.Range("D4").Value = 89.483117936170359
$ws.Range("D5").Value = 68.11196707133476
$ws.Range("D6").Value = 23.123170075151371
$ws.Range("D7").Value = 22.711324889354078
$ws.Range("D8").Value = 41.45269930945453
$ws.Range("D9").Value = 31.806834540952146
$ws.Range("D10").Value = 40.12449432800058
$ws.Range("D11").Value = 36.320819037737238
$ws.Range("D12").Value = 25.704092090099973
$ws.Range("D13").Value = 8.4398132210442682
$ws.Range("D14").Value = 10.516882063829636
$ws.Range("D15").Value = 31.888032928665233
$ws.Range("D16").Value = 76.876829924848636
$ws.Range("D17").Value = 77.288675110645926
$ws.Range("D18").Value = 58.54730069054547
$ws.Range("D19").Value = 68.193165459047862
$ws.Range("D20").Value = 59.87550567199942
$ws.Range("D21").Value = 63.679180962262762

# Widen column C (now showing longer numeric values) to fit its contents
$ws.Columns.Item(3).ColumnWidth = 17.7

# Leave the selection where the author last clicked while building the chart
$ws.Range("L12").Select()
